$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.511.39'
$ws.Range("E2").Value = '  -2.92%  '

$ws.Range("D3").Value = '3.375.08'
$ws.Range("E3").Value = '  -3.85%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''571.30'
$ws.Range("E5").Value = '  -3.48%  '

$ws.Range("D6").Value = '''125.33'
$ws.Range("E6").Value = '  -6.98%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '3.374.99'
$ws.Range("E8").Value = '  -3.81%  '

$ws.Range("E9").Value = '  -2.65%  '

$ws.Range("D10").Value = '''7.26'
$ws.Range("E10").Value = '  -4.36%  '

$ws.Range("E11").Value = '  -4.99%  '

$ws.Range("D12").Value = '''0.374'
$ws.Range("E12").Value = '  -3.74%  '

$ws.Range("D13").Value = '3.954.43'
$ws.Range("E13").Value = '  -3.73%  '

$ws.Range("E14").Value = '  -0.96%  '

$ws.Range("D15").Value = '3.380.09'
$ws.Range("E15").Value = '  -3.68%  '

$ws.Range("E16").Value = '  -6.32%  '

$ws.Range("D17").Value = '62.580.10'
$ws.Range("E17").Value = '  -2.79%  '

$ws.Range("D18").Value = '''24.35'
$ws.Range("E18").Value = '  -5.49%  '

$ws.Range("D19").Value = '''9.21'
$ws.Range("E19").Value = '  -7.71%  '

$ws.Range("D20").Value = '''5.61'
$ws.Range("E20").Value = '  -2.51%  '

$ws.Range("D21").Value = '''13.05'
$ws.Range("E21").Value = '  -4.36%  '

$ws.Range("D22").Value = '''371.53'
$ws.Range("E22").Value = '  -4.56%  '

$ws.Range("D23").Value = '''0.552'
$ws.Range("E23").Value = '  -4.78%  '

$ws.Range("D24").Value = '3.511.89'
$ws.Range("E24").Value = '  -3.78%  '

$ws.Range("D25").Value = '''0.999'
$ws.Range("E25").Value = '  -0.16%  '

$ws.Range("D26").Value = '''71.61'
$ws.Range("E26").Value = '  -3.97%  '

$ws.Range("D27").Value = '''0.0000105'
$ws.Range("E27").Value = '  -10.66%  '

$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.11%  '

$ws.Range("D29").Value = '''6.88'
$ws.Range("E29").Value = '  -7.30%  '

$ws.Range("D30").Value = '''2.10'
$ws.Range("E30").Value = '  -7.44%  '

$ws.Range("D31").Value = '''7.76'
$ws.Range("E31").Value = '  -6.14%  '

$ws.Range("B33").Value = 'RenzoRestakedETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D33").Value = '3.409.81'
$ws.Range("E33").Value = '  -3.57%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '''0.148'
$ws.Range("E34").Value = '  -5.45%  '

$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '''1.36'
$ws.Range("E35").Value = '  -7.27%  '

$ws.Range("D36").Value = '''22.63'
$ws.Range("E36").Value = '  -3.14%  '

$ws.Range("E37").Value = '  -1.52%  '

$ws.Range("D38").Value = '''166.53'
$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("D39").Value = '''6.60'
$ws.Range("E39").Value = '  -5.21%  '

$ws.Range("E40").Value = '  -5.78%  '

$ws.Range("D41").Value = '''0.0751'
$ws.Range("E41").Value = '  -4.70%  '

$ws.Range("E42").Value = '  +0.07%  '

$ws.Range("D43").Value = '''0.763'
$ws.Range("E43").Value = '  -5.94%  '

$ws.Range("D44").Value = '''41.47'
$ws.Range("E44").Value = '  -1.31%  '

$ws.Range("D45").Value = '''4.22'
$ws.Range("E45").Value = '  -5.28%  '

$ws.Range("D46").Value = '''22.51'
$ws.Range("E46").Value = '  -9.40%  '

$ws.Range("D47").Value = '''1.53'
$ws.Range("E47").Value = '  -7.72%  '

$ws.Range("D48").Value = '''1.07'
$ws.Range("E48").Value = '  -8.47%  '

$ws.Range("D49").Value = '''6.57'
$ws.Range("E49").Value = '  -3.36%  '

$ws.Range("D50").Value = '2.236.16'
$ws.Range("E50").Value = '  -5.76%  '

$ws.Range("D51").Value = '''0.837'
$ws.Range("E51").Value = '  -7.93%  '

# Reset style to Normal for text-forced numeric-looking cells (removes quotePrefix style index)
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
